$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 'last updated' timestamp in A1 (covid stats refreshed to 02:05)
$ws.Range("A1").Value = 'Datos actualizados a 1 de Junio de 2020 a las 02:05'

# Country data table: Row, Pais, CasosTotales, NuevosCasos, CasosActivos, Recuperados, CasosCriticos, MuertesHoy, Muertes
$data = @(
  @(4, 'Estados Unidos', 1837165, 20345, 541388, 1189587, 0, 633, 106190),
  @(5, 'Brasil', 514849, 16409, 206555, 278980, 0, 480, 29314),
  @(6, 'Rusia', 405843, 9268, 171883, 229267, 0, 138, 4693),
  @(7, 'España', 286509, 201, 196958, 62424, 0, 2, 27127),
  @(8, 'Reino Unido', 274762, 1936, 0, 0, 0, 113, 38489),
  @(9, 'Italia', 232997, 333, 157507, 42075, 0, 75, 33415),
  @(10, 'India', 190609, 8782, 91852, 93349, 0, 223, 5408),
  @(11, 'Francia', 188882, 257, 68355, 91725, 0, 31, 28802),
  @(12, 'Alemania', 183494, 200, 165200, 9689, 0, 5, 8605),
  @(13, 'Peru', 164476, 8805, 67208, 92762, 0, 135, 4506),
  @(14, 'Turquia', 163942, 839, 127973, 31429, 0, 25, 4540),
  @(15, 'Iran', 151466, 2516, 118848, 24821, 0, 63, 7797),
  @(16, 'Chile', 99688, 4830, 42727, 55907, 0, 57, 1054),
  @(17, 'Canada', 90928, 738, 48839, 34795, 0, 221, 7294),
  @(18, 'Mexico', 87512, 2885, 61871, 15862, 0, 364, 9779),
  @(19, 'Arabia Saudita', 85261, 1877, 62442, 22316, 0, 23, 503),
  @(20, 'China', 83001, 2, 78304, 63, 0, 0, 4634),
  @(21, 'Pakistan', 69496, 3039, 25271, 42742, 0, 88, 1483),
  @(22, 'Belgica', 58381, 195, 15887, 33027, 0, 14, 9467),
  @(23, 'Catar', 56910, 1648, 30290, 26582, 0, 2, 38),
  @(24, 'Banglades', 47153, 2545, 9781, 36722, 0, 40, 650),
  @(25, 'Paises Bajos', 46442, 185, 0, 0, 0, 5, 5956),
  @(26, 'Bielorrusia', 42556, 898, 18514, 23807, 0, 6, 235),
  @(27, 'Ecuador', 39098, 527, 19592, 16148, 0, 24, 3358),
  @(28, 'Suecia', 37542, 429, 4971, 28176, 0, 0, 4395),
  @(29, 'Singapur', 34884, 518, 21699, 13162, 0, 0, 23),
  @(30, 'Emiratos Arabes Unidos', 34557, 661, 17932, 16361, 0, 2, 264),
  @(31, 'Sudafrica', 32683, 1716, 16809, 15191, 0, 40, 683),
  @(32, 'Portugal', 32500, 297, 19409, 11681, 0, 14, 1410),
  @(33, 'Suiza', 30862, 17, 28500, 442, 0, 1, 1920),
  @(34, 'Colombia', 29383, 1147, 8543, 19901, 0, 49, 939),
  @(35, 'Kuwait', 27043, 851, 11386, 15445, 0, 7, 212),
  @(36, 'Indonesia', 26473, 700, 7308, 17552, 0, 40, 1613),
  @(37, 'Irlanda', 24990, 61, 22089, 1249, 0, 1, 1652),
  @(38, 'Egipto', 24985, 1536, 6037, 17989, 0, 46, 959),
  @(39, 'Polonia', 23786, 215, 11271, 11451, 0, 3, 1064),
  @(40, 'Ucrania', 23672, 468, 9538, 13426, 0, 12, 708),
  @(41, 'Rumania', 19257, 124, 13256, 4735, 0, 7, 1266),
  @(42, 'Filipinas', 18086, 862, 3909, 13220, 0, 7, 957),
  @(43, 'Republica Dominicana', 17285, 377, 10559, 6224, 0, 4, 502),
  @(44, 'Israel', 17071, 59, 14812, 1974, 0, 1, 285),
  @(45, 'Argentina', 16851, 637, 5336, 10976, 0, 11, 539),
  @(46, 'Japon', 16851, 47, 14459, 1501, 0, 5, 891),
  @(47, 'Austria', 16731, 46, 15593, 470, 0, 0, 668),
  @(48, 'Afganistan', 15205, 680, 1328, 13620, 0, 8, 257),
  @(49, 'Panama', 13463, 445, 9514, 3613, 0, 6, 336),
  @(50, 'Dinamarca', 11669, 36, 10362, 733, 0, 3, 574),
  @(51, 'Corea del Sur', 11468, 27, 10405, 793, 0, 1, 270),
  @(52, 'Oman', 11437, 1014, 2682, 8706, 0, 7, 49),
  @(53, 'Serbia', 11412, 31, 6698, 4471, 0, 1, 243),
  @(54, 'Barein', 11398, 605, 6673, 4706, 0, 2, 19),
  @(55, 'Kazajistan', 10858, 476, 5404, 5414, 0, 2, 40),
  @(56, 'Nigeria', 10162, 307, 3007, 6868, 0, 14, 287),
  @(57, 'Bolivia', 9592, 861, 889, 8393, 0, 10, 310),
  @(58, 'Argelia', 9394, 127, 5748, 2993, 0, 7, 653),
  @(59, 'Armenia', 9282, 355, 3386, 5765, 0, 4, 131),
  @(60, 'Chequia', 9268, 38, 6558, 2390, 0, 1, 320),
  @(61, 'Noruega', 8440, 3, 7727, 477, 0, 0, 236),
  @(62, 'Moldavia', 8251, 153, 4581, 3375, 0, 4, 295),
  @(63, 'Ghana', 7881, 113, 2841, 5004, 0, 1, 36),
  @(64, 'Malasia', 7819, 57, 6353, 1351, 0, 0, 115),
  @(65, 'Marruecos', 7807, 27, 5459, 2143, 0, 1, 205),
  @(66, 'Australia', 7195, 10, 6614, 478, 0, 0, 103),
  @(67, 'Finlandia', 6859, 33, 5500, 1039, 0, 4, 320),
  @(68, 'Irak', 6439, 260, 3156, 3078, 0, 10, 205),
  @(69, 'Camerun', 5904, 0, 3568, 2145, 0, 0, 191),
  @(70, 'Azerbaiyan', 5494, 248, 3428, 2003, 0, 2, 63),
  @(71, 'Honduras', 5094, 208, 536, 4357, 0, 2, 201),
  @(72, 'Sudan', 5026, 226, 1423, 3317, 0, 24, 286),
  @(73, 'Guatemala', 4739, 132, 706, 3931, 0, 12, 102),
  @(74, 'Luxemburgo', 4018, 2, 3833, 75, 0, 0, 110),
  @(75, 'Tayikistan', 3930, 123, 2004, 1879, 0, 0, 47),
  @(76, 'Hungria', 3876, 9, 2147, 1203, 0, 2, 526),
  @(77, 'Guinea', 3706, 0, 2030, 1653, 0, 0, 23),
  @(78, 'Senegal', 3645, 110, 1801, 1802, 0, 0, 42),
  @(79, 'Uzbekistan', 3623, 77, 2837, 771, 0, 1, 15),
  @(80, 'Republica de Yibuti', 3354, 160, 1504, 1826, 0, 2, 24),
  @(81, 'Tailandia', 3081, 4, 2963, 61, 0, 0, 57),
  @(82, 'Consejo Danes para los Refugiados', 3070, 104, 448, 2550, 0, 3, 72),
  @(83, 'Grecia', 2917, 2, 1374, 1368, 0, 0, 175),
  @(84, 'Costa de Marfil', 2833, 34, 1435, 1365, 0, 0, 33),
  @(85, 'Gabon', 2655, 0, 722, 1916, 0, 0, 17),
  @(86, 'El Salvador', 2517, 122, 1040, 1431, 0, 2, 46),
  @(87, 'Bulgaria', 2513, 14, 1074, 1299, 0, 1, 140),
  @(88, 'Bosnia y Herzegovina', 2510, 16, 1862, 495, 0, 0, 153),
  @(89, 'Croacia', 2246, 0, 2072, 71, 0, 0, 103),
  @(90, 'Republica de Macedonia', 2226, 62, 1552, 541, 0, 2, 133),
  @(91, 'Cuba', 2045, 20, 1809, 153, 0, 0, 83),
  @(92, 'Somalia', 1976, 60, 348, 1550, 0, 5, 78),
  @(93, 'Kenia', 1962, 74, 478, 1420, 0, 1, 64),
  @(94, 'Estonia', 1869, 4, 1624, 177, 0, 1, 68),
  @(95, 'Haiti', 1865, 281, 24, 1800, 0, 6, 41),
  @(96, 'Islandia', 1806, 0, 1794, 2, 0, 0, 10),
  @(97, 'Maldivas', 1773, 101, 453, 1315, 0, 0, 5),
  @(98, 'Kirguistan', 1748, 26, 1170, 562, 0, 0, 16),
  @(99, 'Mayotte', 1699, 0, 1385, 293, 0, 0, 21),
  @(100, 'Lituania', 1675, 5, 1236, 369, 0, 0, 70),
  @(101, 'Sri Lanka', 1633, 20, 801, 822, 0, 0, 10),
  @(102, 'Nepal', 1572, 171, 220, 1344, 0, 2, 8),
  @(103, 'Eslovaquia', 1521, 0, 1366, 127, 0, 0, 28),
  @(104, 'Venezuela', 1510, 51, 302, 1194, 0, 0, 14),
  @(105, 'Nueva Zelanda', 1504, 0, 1481, 1, 0, 0, 22),
  @(106, 'Eslovenia', 1473, 0, 1358, 7, 0, 0, 108),
  @(107, 'Guinea Ecuatorial', 1306, 0, 200, 1094, 0, 0, 12),
  @(108, 'Mali', 1265, 15, 716, 472, 0, 1, 77),
  @(109, 'Guinea-Bisau', 1256, 0, 42, 1206, 0, 0, 8),
  @(110, 'Libano', 1220, 29, 712, 481, 0, 1, 27),
  @(111, 'Etiopia', 1172, 109, 209, 952, 0, 3, 11),
  @(112, 'Albania', 1137, 15, 872, 232, 0, 0, 33),
  @(113, 'Hong Kong', 1085, 2, 1037, 44, 0, 0, 4),
  @(114, 'Tunez', 1077, 1, 960, 69, 0, 0, 48),
  @(115, 'Letonia', 1066, 1, 745, 297, 0, 0, 24),
  @(116, 'Zambia', 1057, 0, 779, 271, 0, 0, 7),
  @(117, 'Costa Rica', 1056, 9, 669, 377, 0, 0, 10),
  @(118, 'Republica de Africa Central', 1011, 49, 23, 986, 0, 1, 2),
  @(119, 'Sudan del Sur', 994, 0, 6, 978, 0, 0, 10),
  @(120, 'Paraguay', 986, 22, 477, 498, 0, 0, 11),
  @(121, 'Niger', 958, 2, 839, 55, 0, 0, 64),
  @(122, 'Republica de Chipre', 944, 1, 790, 137, 0, 0, 17),
  @(123, 'Sierra Leona', 861, 9, 454, 361, 0, 0, 46),
  @(124, 'Burkina Faso', 847, 0, 720, 74, 0, 0, 53),
  @(125, 'Uruguay', 823, 2, 685, 116, 0, 0, 22),
  @(126, 'Georgia', 783, 26, 605, 166, 0, 0, 12),
  @(127, 'Republica del Chad', 778, 19, 491, 222, 0, 0, 65),
  @(128, 'Madagascar', 771, 13, 168, 597, 0, 0, 6),
  @(129, 'Principado de Andorra', 764, 0, 694, 19, 0, 0, 51),
  @(130, 'Nicaragua', 759, 0, 370, 354, 0, 0, 35),
  @(131, 'Jordania', 739, 5, 522, 208, 0, 0, 9),
  @(132, 'Crucero', 712, 0, 651, 48, 0, 0, 13),
  @(133, 'San Marino', 671, 0, 357, 272, 0, 0, 42),
  @(134, 'Malta', 618, 0, 534, 75, 0, 0, 9),
  @(135, 'Congo', 611, 40, 179, 412, 0, 1, 20),
  @(136, 'Jamaica', 581, 6, 290, 282, 0, 0, 9),
  @(137, 'Mauritania', 530, 47, 27, 480, 0, 3, 23),
  @(138, 'Tanzania', 509, 0, 183, 305, 0, 0, 21),
  @(139, 'Guayana Francesa', 499, 22, 200, 298, 0, 0, 1),
  @(140, 'Santo Tome y Principe', 483, 4, 68, 403, 0, 0, 12),
  @(141, 'Reunion', 471, 0, 411, 59, 0, 0, 1),
  @(142, 'Estado de Palestina', 448, 1, 372, 73, 0, 0, 3),
  @(143, 'Togo', 442, 9, 211, 218, 0, 0, 13),
  @(144, 'Taiwan', 442, 0, 423, 12, 0, 0, 7),
  @(145, 'Cabo Verde', 435, 14, 193, 238, 0, 0, 4),
  @(146, 'Uganda', 417, 4, 72, 345, 0, 0, 0),
  @(147, 'Ruanda', 370, 11, 256, 113, 0, 0, 1),
  @(148, 'Isla de Man', 336, 0, 309, 3, 0, 0, 24),
  @(149, 'Mauricio', 335, 0, 322, 3, 0, 0, 10),
  @(150, 'Vietnam', 328, 0, 279, 49, 0, 0, 0),
  @(151, 'Montenegro', 324, 0, 315, 0, 0, 0, 9),
  @(152, 'Yemen', 323, 13, 14, 229, 0, 3, 80),
  @(153, 'Liberia', 288, 8, 157, 104, 0, 0, 27),
  @(154, 'Suazilandia', 285, 2, 189, 94, 0, 0, 2),
  @(155, 'Malaui', 284, 5, 42, 238, 0, 0, 4),
  @(156, 'Mozambique', 254, 10, 91, 161, 0, 0, 2),
  @(157, 'Benin', 232, 8, 143, 86, 0, 0, 3),
  @(158, 'Birmania', 224, 0, 138, 80, 0, 0, 6),
  @(159, 'Martinica', 200, 0, 98, 88, 0, 0, 14),
  @(160, 'Islas Feroe', 187, 0, 187, 0, 0, 0, 0),
  @(161, 'Mongolia', 179, 0, 44, 135, 0, 0, 0),
  @(162, 'Zimbabue', 178, 4, 29, 145, 0, 0, 4),
  @(163, 'Gibraltar', 170, 1, 149, 21, 0, 0, 0),
  @(164, 'Guadalupe', 162, 0, 138, 10, 0, 0, 14),
  @(165, 'Libia', 156, 26, 52, 99, 0, 0, 5),
  @(166, 'Guyana', 152, 0, 67, 73, 0, 0, 12),
  @(167, 'Islas Caimanes', 141, 0, 68, 72, 0, 0, 1),
  @(168, 'Brunei', 141, 0, 138, 1, 0, 0, 2),
  @(169, 'Bermudas', 140, 0, 92, 39, 0, 0, 9),
  @(170, 'Camboya', 125, 0, 123, 2, 0, 0, 0),
  @(171, 'Siria', 122, 0, 46, 71, 0, 1, 5),
  @(172, 'Trinidad yTobago', 117, 1, 108, 1, 0, 0, 8),
  @(173, 'Comoras', 106, 0, 26, 78, 0, 0, 2),
  @(174, 'Bahamas', 102, 0, 48, 43, 0, 0, 11),
  @(175, 'Aruba', 101, 0, 98, 0, 0, 0, 3),
  @(176, 'Monaco', 99, 0, 90, 5, 0, 0, 4),
  @(177, 'Barbados', 92, 0, 76, 9, 0, 0, 7),
  @(178, 'Angola', 86, 2, 18, 64, 0, 0, 4),
  @(179, 'Liechtenstein', 82, 0, 55, 26, 0, 0, 1),
  @(180, 'San Martin (Parte Holandesa)', 77, 0, 60, 2, 0, 0, 15),
  @(181, 'Burundi', 63, 0, 33, 29, 0, 0, 1),
  @(182, 'Polinesia Francesa', 60, 0, 60, 0, 0, 0, 0),
  @(183, 'Macao', 45, 0, 45, 0, 0, 0, 0),
  @(184, 'Butan', 43, 10, 6, 37, 0, 0, 0),
  @(185, 'San Martin (Parte Francesa)', 41, 0, 33, 5, 0, 0, 3),
  @(186, 'Puerto Rico', 39, 0, 1, 36, 0, 0, 2),
  @(187, 'Eritrea', 39, 0, 39, 0, 0, 0, 0),
  @(188, 'Botsuana', 35, 0, 20, 14, 0, 0, 1),
  @(189, 'Guam', 32, 0, 0, 31, 0, 0, 1),
  @(190, 'San Vicente y las Granadinas', 26, 0, 15, 11, 0, 0, 0),
  @(191, 'Antigua y Barbuda', 26, 1, 19, 4, 0, 0, 3),
  @(192, 'Gambia', 25, 0, 20, 4, 0, 0, 1),
  @(193, 'Namibia', 24, 1, 14, 10, 0, 0, 0),
  @(194, 'Timor Oriental', 24, 0, 24, 0, 0, 0, 0),
  @(195, 'Surinam', 23, 11, 9, 13, 0, 0, 1),
  @(196, 'Granada', 23, 0, 18, 5, 0, 0, 0),
  @(197, 'Curazao', 19, 0, 14, 4, 0, 0, 1),
  @(198, 'Laos', 19, 0, 16, 3, 0, 0, 0),
  @(199, 'Nueva Caledonia', 19, 0, 18, 1, 0, 0, 0),
  @(200, 'Fiyi', 18, 0, 15, 3, 0, 0, 0),
  @(201, 'Santa Lucia', 18, 0, 18, 0, 0, 0, 0),
  @(202, 'Belice', 18, 0, 16, 0, 0, 0, 2),
  @(203, 'Islas Virgenes de los Estados Unidos', 17, 0, 0, 17, 0, 0, 0),
  @(204, 'Dominica', 16, 0, 16, 0, 0, 0, 0),
  @(205, 'San Cristobal y Nieves', 15, 0, 15, 0, 0, 0, 0),
  @(206, 'Groenlandia', 13, 0, 11, 2, 0, 0, 0),
  @(207, 'Islas Malvinas', 13, 0, 13, 0, 0, 0, 0),
  @(208, 'Santa Sede', 12, 0, 2, 10, 0, 0, 0),
  @(209, 'Islas Turcas y Caicos', 12, 0, 11, 0, 0, 0, 1),
  @(210, 'Montserrat', 11, 0, 10, 0, 0, 0, 1),
  @(211, 'Seychelles', 11, 0, 11, 0, 0, 0, 0),
  @(212, 'Sahara Occidental', 9, 0, 6, 2, 0, 0, 1),
  @(213, 'Islas Virgenes Britanicas', 8, 0, 7, 0, 0, 0, 1),
  @(214, 'Papua Nueva Guinea', 8, 0, 8, 0, 0, 0, 0),
  @(215, 'San Bartolome', 6, 0, 6, 0, 0, 0, 0),
  @(216, 'Bonaire, San Eustaquio y Saba', 6, 0, 6, 0, 0, 0, 0),
  @(217, 'Anguila', 3, 0, 3, 0, 0, 0, 0),
  @(218, 'Lesoto', 2, 0, 1, 1, 0, 0, 0),
  @(219, 'San Pedro y Miquelon', 1, 0, 1, 0, 0, 0, 0)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
    $ws.Cells.Item($r, 5).Value = $item[5]
    $ws.Cells.Item($r, 6).Value = $item[6]
    $ws.Cells.Item($r, 7).Value = $item[7]
    $ws.Cells.Item($r, 8).Value = $item[8]
}
